$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Adder-BOM")

# Set "Bought" (column D) quantities for a handful of parts - BOM update
$ws.Range("D3").Value = 20
$ws.Range("D12").Value = 10
$ws.Range("D15").Value = 2
$ws.Range("D16").Value = 2

# Move selection / scroll position of the frozen pane
$ws.Range("A2").Select()
$excel.ActiveWindow.ScrollColumn = 7
